$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for all data rows (2-10): 46062 -> 46063
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 46063
}

# Capture current (before) values for the rows whose contents get
# swapped/rotated: A (Beteckning), B (Datum), G (Area (ha))
# NOTE: use Value2 (not Value) to read - Value getter is unreliable here.
$row4 = @($ws.Cells.Item(4, 1).Value2, $ws.Cells.Item(4, 2).Value2, $ws.Cells.Item(4, 7).Value2)
$row5 = @($ws.Cells.Item(5, 1).Value2, $ws.Cells.Item(5, 2).Value2, $ws.Cells.Item(5, 7).Value2)
$row7 = @($ws.Cells.Item(7, 1).Value2, $ws.Cells.Item(7, 2).Value2, $ws.Cells.Item(7, 7).Value2)
$row8 = @($ws.Cells.Item(8, 1).Value2, $ws.Cells.Item(8, 2).Value2, $ws.Cells.Item(8, 7).Value2)
$row9 = @($ws.Cells.Item(9, 1).Value2, $ws.Cells.Item(9, 2).Value2, $ws.Cells.Item(9, 7).Value2)

# Rows 4 and 7 swap their contents
$ws.Cells.Item(7, 1).Value = $row4[0]
$ws.Cells.Item(7, 2).Value = $row4[1]
$ws.Cells.Item(7, 7).Value = $row4[2]

$ws.Cells.Item(4, 1).Value = $row7[0]
$ws.Cells.Item(4, 2).Value = $row7[1]
$ws.Cells.Item(4, 7).Value = $row7[2]

# Rows 5, 8, 9 rotate: 5 -> 8, 8 -> 9, 9 -> 5
$ws.Cells.Item(8, 1).Value = $row5[0]
$ws.Cells.Item(8, 2).Value = $row5[1]
$ws.Cells.Item(8, 7).Value = $row5[2]

$ws.Cells.Item(9, 1).Value = $row8[0]
$ws.Cells.Item(9, 2).Value = $row8[1]
$ws.Cells.Item(9, 7).Value = $row8[2]

$ws.Cells.Item(5, 1).Value = $row9[0]
$ws.Cells.Item(5, 2).Value = $row9[1]
$ws.Cells.Item(5, 7).Value = $row9[2]
